$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column for rows 2-5
# from serial 45233 (2023-11-03) to serial 45243 (2023-11-13)
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 3).Value = 45243
}
